# "Pour la prochaine fois" sheet - update the notes list:
#  - Remove outdated items about enemy-above-box, perso animation, enemy entity handlers
#  - Add new findings about the enemy hitbox / kinematic controller bug
#  - Fix the typo "spcler" -> "socket"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the previous content entirely, then rewrite the list at its new
# row positions.
$ws.Range("B2:C18").ClearContents()

$ws.Range("B3").Value = "Joueur enemi sur le sol OK, mais hit box tjrs en l'air"
$ws.Range("C3").Value = "pb : si on la baisse, elle bug avec le sol et tous les joueurs bug (kinematic controler)"

$ws.Range("B7").Value = "deleteRemotePlayer : enlever le perso de la scene"

$ws.Range("B9").Value = "back/utils/socket : connection => prend le dernier user : faire avec pseudo"

$ws.Range("B11").Value = "afficher pseudo sur skins"

$ws.Range("B13").Value = "POV d'un joueur bug desfois (mouvement souris trop brusque ?)"

$ws.Range("B15").Value = "optimiser le code en general"
$ws.Range("C15").Value = "pour une prochiane co, sauvegarder dans le temps"

$ws.Range("B17").Value = "Sauvegarder kd par joueurs ? "

$ws.Range("B19").Value = "Faille xss/autre securités"

$ws.Range("B5").Select() | Out-Null
